# Update the "Learn Status" table from a 9-column per-skill breakdown down
# to a compact 3-column (Name | Not on Learn | Total) report with date-range
# style rollups.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shp = $s.Shapes.Item(3)
$tbl = $shp.Table

# 1. Drop the six per-skill columns (Communication Skills .. Listening
#    Skills), keeping only Name, "Not on Learn", and "Total".
for ($i = 0; $i -lt 6; $i++) {
    $tbl.Columns.Item(2).Delete()
}

# 2. Re-balance the remaining three columns to fill the same overall
#    table width (previously 9 x 76pt, now 3 x 228pt == 8686800 EMU).
for ($c = 1; $c -le $tbl.Columns.Count; $c++) {
    $tbl.Columns.Item($c).Width = 228
}

# At this point the rows read:
#   1: ""              | Not on Learn | Total
#   2: Hyderabad        | 5            | 39
#   3: Bangalore        | 16425        | 16625
#   4: Administration   | 2            | 2
#   5: Delhi            | 4            | 4
#   6: Total            | 16436        | 16670

# 3. Remove the old granular Bangalore row...
$tbl.Rows.Item(3).Delete()

# 4. ...and reinsert it ahead of Hyderabad with its rolled-up totals.
$tbl.Rows.Add(2)
$tbl.Cell(2, 1).Shape.TextFrame.TextRange.Text = "Bangalore"
$tbl.Cell(2, 2).Shape.TextFrame.TextRange.Text = "16625"
$tbl.Cell(2, 3).Shape.TextFrame.TextRange.Text = "16625"

# 5. Hyderabad's "Not on Learn" figure becomes its (now single) total.
$tbl.Cell(3, 2).Shape.TextFrame.TextRange.Text = "39"

# 6. The grand total's "Not on Learn" figure likewise becomes the new total.
$tbl.Cell(6, 2).Shape.TextFrame.TextRange.Text = "16670"
